$d = $word.ActiveDocument

$d.Content.Find.Execute("carte che i giocatori hanno in mano", $true, $false, $false, $false, $false,
                         $true, 1, $false, "carte che ogni giocatore ha in mano inizialmente", 2)
